# "adding team gradebook as dp2 makeup"
# Slide 1 ("CSSE 220 / Intro to Java Graphics") has a text box named
# "TextBox 1" that shows the day's attendance password. Update the
# password from "howtodraw" to "swing", keeping the existing
# formatting (size/highlight) on that run untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$passwordShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 1") {
        $passwordShape = $shp
        break
    }
}

$tr = $passwordShape.TextFrame.TextRange
$passwordParagraph = $tr.Paragraphs(2)
$passwordParagraph.Runs(1).Text = "swing"
